$wb = $excel.ActiveWorkbook

# zh-cn sheet: rows 3 and 4 get updated Handoff/Handback datetimes
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E3").Value = "2016-03-13 02:18:49"
$wsZh.Range("H3").Value = "2016-03-13 02:19:08"
$wsZh.Range("E4").Value = "2016-03-13 02:18:49"
$wsZh.Range("H4").Value = "2016-03-13 02:19:08"

# de-de sheet: rows 3 and 4 get updated Handoff/Handback datetimes
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E3").Value = "2016-03-13 02:18:52"
$wsDe.Range("H3").Value = "2016-03-13 02:19:14"
$wsDe.Range("E4").Value = "2016-03-13 02:18:52"
$wsDe.Range("H4").Value = "2016-03-13 02:19:14"
